$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 23.809525  # H2: 26.666666 -> 23.809525
$ws.Cells.Item(2, 9).Value = 24  # I2: 27.272728 -> 24
$ws.Cells.Item(2, 11).Value = 24  # K2: 27.272728 -> 24
$ws.Cells.Item(2, 13).Value = 89  # M2: 85.727272 -> 89
$ws.Cells.Item(28, 8).Value = 897.7143  # H28: 858 -> 897.7143
$ws.Cells.Item(28, 9).Value = 565.7895  # I28: 482.35294 -> 565.7895
$ws.Cells.Item(28, 11).Value = 565.7895  # K28: 482.35294 -> 565.7895
$ws.Cells.Item(28, 13).Value = -80.78949999999998  # M28: 2.64706000000001 -> -80.78949999999998
$ws.Cells.Item(87, 8).Value = 54959  # H87: 58087.89 -> 54959
$ws.Cells.Item(87, 10).Value = 54959  # J87: 58087.89 -> 54959
$ws.Cells.Item(87, 12).Value = 54959  # L87: 58087.89 -> 54959
$ws.Cells.Item(87, 14).Value = -57455  # N87: -60583.89 -> -57455
$ws.Cells.Item(90, 8).Value = 54959  # H90: 58087.89 -> 54959
$ws.Cells.Item(90, 10).Value = 54959  # J90: 58087.89 -> 54959
$ws.Cells.Item(90, 12).Value = 164877  # L90: 174263.67 -> 164877
$ws.Cells.Item(90, 14).Value = -177357  # N90: -186743.67 -> -177357
$ws.Cells.Item(100, 8).Value = 3226.9167  # H100: 2856.2856 -> 3226.9167
$ws.Cells.Item(100, 9).Value = 3174.25  # I100: 2899.4 -> 3174.25
$ws.Cells.Item(100, 10).Value = 3332.25  # J100: 2748.5 -> 3332.25
$ws.Cells.Item(100, 11).Value = 3174.25  # K100: 2899.4 -> 3174.25
$ws.Cells.Item(100, 12).Value = 3332.25  # L100: 2748.5 -> 3332.25
$ws.Cells.Item(100, 13).Value = -2633.25  # M100: -2358.4 -> -2633.25
$ws.Cells.Item(100, 14).Value = -4414.25  # N100: -3830.5 -> -4414.25
$ws.Cells.Item(112, 8).Value = 4062.5  # H112: 4260.0356 -> 4062.5
$ws.Cells.Item(112, 10).Value = 4395.407  # J112: 4643.28 -> 4395.407
$ws.Cells.Item(112, 12).Value = 13186.221  # L112: 13929.84 -> 13186.221
$ws.Cells.Item(112, 14).Value = -15402.221  # N112: -16145.84 -> -15402.221
$ws.Cells.Item(117, 8).Value = 14800  # H117: 14900 -> 14800
$ws.Cells.Item(117, 10).Value = 14800  # J117: 14900 -> 14800
$ws.Cells.Item(117, 12).Value = 14800  # L117: 14900 -> 14800
$ws.Cells.Item(117, 14).Value = -23978  # N117: -24078 -> -23978
$ws.Cells.Item(132, 8).Value = 1036.4642  # H132: 1003.62067 -> 1036.4642
$ws.Cells.Item(132, 9).Value = 963.7406999999999  # I132: 932.3214 -> 963.7406999999999
$ws.Cells.Item(132, 11).Value = 2891.2221  # K132: 2796.9642 -> 2891.2221
$ws.Cells.Item(132, 13).Value = -361.2221  # M132: -266.9642000000003 -> -361.2221

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2252322  # H32: 2277915.8 -> 2252322
$ws.Cells.Item(32, 9).Value = 2277348.2  # I32: 2303524 -> 2277348.2
$ws.Cells.Item(32, 11).Value = 2277348.2  # K32: 2303524 -> 2277348.2
$ws.Cells.Item(32, 13).Value = -2277061.2  # M32: -2303237 -> -2277061.2
$ws.Cells.Item(45, 8).Value = 3694.3462  # H45: 3695.423 -> 3694.3462
$ws.Cells.Item(45, 9).Value = 1875.2858  # I45: 1877.2858 -> 1875.2858
$ws.Cells.Item(45, 11).Value = 1875.2858  # K45: 1877.2858 -> 1875.2858
$ws.Cells.Item(45, 13).Value = -1498.2858  # M45: -1500.2858 -> -1498.2858
$ws.Cells.Item(61, 8).Value = 18871262  # H61: 18871282 -> 18871262
$ws.Cells.Item(61, 9).Value = 2641.6829  # I61: 2690.25 -> 2641.6829
$ws.Cells.Item(61, 10).Value = 83339050  # J61: 76928490 -> 83339050
$ws.Cells.Item(61, 11).Value = 2641.6829  # K61: 2690.25 -> 2641.6829
$ws.Cells.Item(61, 12).Value = 83339050  # L61: 76928490 -> 83339050
$ws.Cells.Item(61, 13).Value = -2429.6829  # M61: -2478.25 -> -2429.6829
$ws.Cells.Item(61, 14).Value = -83339474  # N61: -76928914 -> -83339474
$ws.Cells.Item(136, 8).Value = 18871262  # H136: 18871282 -> 18871262
$ws.Cells.Item(136, 9).Value = 2641.6829  # I136: 2690.25 -> 2641.6829
$ws.Cells.Item(136, 10).Value = 83339050  # J136: 76928490 -> 83339050
$ws.Cells.Item(136, 11).Value = 7925.048699999999  # K136: 8070.75 -> 7925.048699999999
$ws.Cells.Item(136, 12).Value = 250017150  # L136: 230785470 -> 250017150
$ws.Cells.Item(136, 13).Value = -5375.048699999999  # M136: -5520.75 -> -5375.048699999999
$ws.Cells.Item(136, 14).Value = -250022250  # N136: -230790570 -> -250022250

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5659.095  # H31: 5668.595 -> 5659.095
$ws.Cells.Item(31, 9).Value = 2183.5386  # I31: 2214.2307 -> 2183.5386
$ws.Cells.Item(31, 11).Value = 2183.5386  # K31: 2214.2307 -> 2183.5386
$ws.Cells.Item(31, 13).Value = -1888.5386  # M31: -1919.2307 -> -1888.5386
$ws.Cells.Item(34, 8).Value = 5659.095  # H34: 5668.595 -> 5659.095
$ws.Cells.Item(34, 9).Value = 2183.5386  # I34: 2214.2307 -> 2183.5386
$ws.Cells.Item(34, 11).Value = 2183.5386  # K34: 2214.2307 -> 2183.5386
$ws.Cells.Item(34, 13).Value = -1981.5386  # M34: -2012.2307 -> -1981.5386
$ws.Cells.Item(41, 8).Value = 35000  # H41: 5000 -> 35000
$ws.Cells.Item(41, 10).Value = 65000  # J41: 0 -> 65000
$ws.Cells.Item(41, 12).Value = 65000  # L41: 0 -> 65000
$ws.Cells.Item(41, 14).Value = -65856  # N41: None -> -65856
$ws.Cells.Item(62, 8).Value = 6560.885  # H62: 6633.9614 -> 6560.885
$ws.Cells.Item(62, 9).Value = 5741.154  # I62: 5911.3335 -> 5741.154
$ws.Cells.Item(62, 10).Value = 7380.615  # J62: 7253.357 -> 7380.615
$ws.Cells.Item(62, 11).Value = 5741.154  # K62: 5911.3335 -> 5741.154
$ws.Cells.Item(62, 12).Value = 7380.615  # L62: 7253.357 -> 7380.615
$ws.Cells.Item(62, 13).Value = -5117.154  # M62: -5287.3335 -> -5117.154
$ws.Cells.Item(62, 14).Value = -8628.615  # N62: -8501.357 -> -8628.615
$ws.Cells.Item(65, 8).Value = 6560.885  # H65: 6633.9614 -> 6560.885
$ws.Cells.Item(65, 9).Value = 5741.154  # I65: 5911.3335 -> 5741.154
$ws.Cells.Item(65, 10).Value = 7380.615  # J65: 7253.357 -> 7380.615
$ws.Cells.Item(65, 11).Value = 28705.77  # K65: 29556.6675 -> 28705.77
$ws.Cells.Item(65, 12).Value = 36903.075  # L65: 36266.785 -> 36903.075
$ws.Cells.Item(65, 13).Value = -25585.77  # M65: -26436.6675 -> -25585.77
$ws.Cells.Item(65, 14).Value = -43143.075  # N65: -42506.785 -> -43143.075
$ws.Cells.Item(99, 8).Value = 5143.5  # H99: 4917.077 -> 5143.5
$ws.Cells.Item(99, 9).Value = 3253.8333  # I99: 3103.2856 -> 3253.8333
$ws.Cells.Item(99, 11).Value = 3253.8333  # K99: 3103.2856 -> 3253.8333
$ws.Cells.Item(99, 13).Value = -1755.8333  # M99: -1605.2856 -> -1755.8333
$ws.Cells.Item(105, 8).Value = 5953596.5  # H105: 5953754 -> 5953596.5
$ws.Cells.Item(105, 9).Value = 7143566  # I105: 7937228 -> 7143566
$ws.Cells.Item(105, 10).Value = 3749.5  # J105: 3331.3333 -> 3749.5
$ws.Cells.Item(105, 11).Value = 7143566  # K105: 7937228 -> 7143566
$ws.Cells.Item(105, 12).Value = 3749.5  # L105: 3331.3333 -> 3749.5
$ws.Cells.Item(105, 13).Value = -7141819  # M105: -7935481 -> -7141819
$ws.Cells.Item(105, 14).Value = -7243.5  # N105: -6825.3333 -> -7243.5
$ws.Cells.Item(123, 8).Value = 84000  # H123: 90000 -> 84000
$ws.Cells.Item(123, 10).Value = 84000  # J123: 90000 -> 84000
$ws.Cells.Item(123, 12).Value = 84000  # L123: 90000 -> 84000
$ws.Cells.Item(123, 14).Value = -93800  # N123: -99800 -> -93800
$ws.Cells.Item(126, 8).Value = 5143.5  # H126: 4917.077 -> 5143.5
$ws.Cells.Item(126, 9).Value = 3253.8333  # I126: 3103.2856 -> 3253.8333
$ws.Cells.Item(126, 11).Value = 9761.499899999999  # K126: 9309.856800000001 -> 9761.499899999999
$ws.Cells.Item(126, 13).Value = -7291.499899999999  # M126: -6839.856800000001 -> -7291.499899999999
$ws.Cells.Item(132, 8).Value = 6820.44  # H132: 6984.8335 -> 6820.44
$ws.Cells.Item(132, 9).Value = 5230  # I132: 5465.5 -> 5230
$ws.Cells.Item(132, 11).Value = 15690  # K132: 16396.5 -> 15690
$ws.Cells.Item(132, 13).Value = -13160  # M132: -13866.5 -> -13160

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 60435268  # H4: 60435270 -> 60435268
$ws.Cells.Item(5, 8).Value = 1219.125  # H5: 1249.5217 -> 1219.125
$ws.Cells.Item(5, 9).Value = 956.7059  # I5: 984 -> 956.7059
$ws.Cells.Item(5, 11).Value = 2870.1177  # K5: 2952 -> 2870.1177
$ws.Cells.Item(5, 13).Value = -2758.1177  # M5: -2840 -> -2758.1177
$ws.Cells.Item(7, 8).Value = 48  # H7: 0 -> 48
$ws.Cells.Item(7, 9).Value = 57.857143  # I7: 0 -> 57.857143
$ws.Cells.Item(7, 10).Value = 13.5  # J7: 0 -> 13.5
$ws.Cells.Item(7, 11).Value = 173.571429  # K7: 0 -> 173.571429
$ws.Cells.Item(7, 12).Value = 40.5  # L7: 0 -> 40.5
$ws.Cells.Item(7, 13).Value = -61.57142899999999  # M7: None -> -61.57142899999999
$ws.Cells.Item(7, 14).Value = -264.5  # N7: None -> -264.5
$ws.Cells.Item(55, 8).Value = 11125645  # H55: 11123423 -> 11125645
$ws.Cells.Item(55, 10).Value = 14302857  # J55: 14300000 -> 14302857
$ws.Cells.Item(55, 12).Value = 42908571  # L55: 42900000 -> 42908571
$ws.Cells.Item(55, 14).Value = -42908925  # N55: -42900354 -> -42908925
$ws.Cells.Item(61, 8).Value = 521.6667  # H61: 467.4 -> 521.6667
$ws.Cells.Item(61, 9).Value = 257.5  # I61: 229 -> 257.5
$ws.Cells.Item(61, 10).Value = 1050  # J61: 825 -> 1050
$ws.Cells.Item(61, 11).Value = 772.5  # K61: 687 -> 772.5
$ws.Cells.Item(61, 12).Value = 3150  # L61: 2475 -> 3150
$ws.Cells.Item(61, 13).Value = -557.5  # M61: -472 -> -557.5
$ws.Cells.Item(61, 14).Value = -3580  # N61: -2905 -> -3580
$ws.Cells.Item(135, 8).Value = 1219.125  # H135: 1249.5217 -> 1219.125
$ws.Cells.Item(135, 9).Value = 956.7059  # I135: 984 -> 956.7059
$ws.Cells.Item(135, 11).Value = 8610.3531  # K135: 8856 -> 8610.3531
$ws.Cells.Item(135, 13).Value = -6075.3531  # M135: -6321 -> -6075.3531

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 87499.664  # H52: 87142.42999999999 -> 87499.664
$ws.Cells.Item(52, 10).Value = 87499.664  # J52: 87142.42999999999 -> 87499.664
$ws.Cells.Item(52, 12).Value = 87499.664  # L52: 87142.42999999999 -> 87499.664
$ws.Cells.Item(52, 14).Value = -88017.664  # N52: -87660.42999999999 -> -88017.664
$ws.Cells.Item(107, 8).Value = 705.2778  # H107: 741.64703 -> 705.2778
$ws.Cells.Item(107, 9).Value = 305.7  # I107: 306.1 -> 305.7
$ws.Cells.Item(107, 10).Value = 1204.75  # J107: 1363.8572 -> 1204.75
$ws.Cells.Item(107, 11).Value = 305.7  # K107: 306.1 -> 305.7
$ws.Cells.Item(107, 12).Value = 1204.75  # L107: 1363.8572 -> 1204.75
$ws.Cells.Item(107, 13).Value = 1614.3  # M107: 1613.9 -> 1614.3
$ws.Cells.Item(107, 14).Value = -5044.75  # N107: -5203.8572 -> -5044.75
$ws.Cells.Item(122, 8).Value = 1693175.9  # H122: 1775701.6 -> 1693175.9
$ws.Cells.Item(122, 9).Value = 1915412.6  # I122: 2021746.8 -> 1915412.6
$ws.Cells.Item(122, 11).Value = 5746237.800000001  # K122: 6065240.4 -> 5746237.800000001
$ws.Cells.Item(122, 13).Value = -5743787.800000001  # M122: -6062790.4 -> -5743787.800000001
$ws.Cells.Item(126, 8).Value = 5250.857  # H126: 5339.385 -> 5250.857
$ws.Cells.Item(126, 9).Value = 4542.3  # I126: 4652.875 -> 4542.3
$ws.Cells.Item(126, 11).Value = 13626.9  # K126: 13958.625 -> 13626.9
$ws.Cells.Item(126, 13).Value = -11156.9  # M126: -11488.625 -> -11156.9

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 16500  # H42: 8000 -> 16500
$ws.Cells.Item(42, 10).Value = 25000  # J42: 0 -> 25000
$ws.Cells.Item(42, 12).Value = 25000  # L42: 0 -> 25000
$ws.Cells.Item(42, 14).Value = -26126  # N42: None -> -26126
$ws.Cells.Item(48, 8).Value = 50000  # H48: 0 -> 50000
$ws.Cells.Item(48, 10).Value = 50000  # J48: 0 -> 50000
$ws.Cells.Item(48, 12).Value = 50000  # L48: 0 -> 50000
$ws.Cells.Item(48, 14).Value = -51322  # N48: None -> -51322
$ws.Cells.Item(49, 8).Value = 16500  # H49: 8000 -> 16500
$ws.Cells.Item(49, 10).Value = 25000  # J49: 0 -> 25000
$ws.Cells.Item(49, 12).Value = 25000  # L49: 0 -> 25000
$ws.Cells.Item(49, 14).Value = -25294  # N49: None -> -25294
$ws.Cells.Item(135, 10).Value = 105000  # J135: 106666.664 -> 105000
$ws.Cells.Item(135, 12).Value = 105000  # L135: 106666.664 -> 105000
$ws.Cells.Item(135, 14).Value = -115140  # N135: -116806.664 -> -115140
$ws.Cells.Item(136, 8).Value = 8724  # H136: 8989.094999999999 -> 8724
$ws.Cells.Item(136, 9).Value = 3536.4614  # I136: 3621.96 -> 3536.4614
$ws.Cells.Item(136, 10).Value = 13374.896  # J136: 13781.179 -> 13374.896
$ws.Cells.Item(136, 11).Value = 10609.3842  # K136: 10865.88 -> 10609.3842
$ws.Cells.Item(136, 12).Value = 40124.688  # L136: 41343.537 -> 40124.688
$ws.Cells.Item(136, 13).Value = -8059.3842  # M136: -8315.880000000001 -> -8059.3842
$ws.Cells.Item(136, 14).Value = -45224.688  # N136: -46443.537 -> -45224.688

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5098.744  # H132: 5293.927 -> 5098.744
$ws.Cells.Item(132, 10).Value = 4708  # J132: 5223.7856 -> 4708
$ws.Cells.Item(132, 12).Value = 14124  # L132: 15671.3568 -> 14124
$ws.Cells.Item(132, 14).Value = -19184  # N132: -20731.3568 -> -19184

Write-Host "Applied all market-data updates."